# Generate Report for Handoff
#
# Refresh the handoff-status report: the "d3763608-e13e-4045-bb9c-c69b6df4aed9.md"
# row got a new xliff generated during this handoff run, so its timestamps
# move forward on both the Overview sheet and the zh-cn detail sheet.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for d3763608-...md (row 6)
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G6").Value = "2016-08-24 20:43:12"

# zh-cn detail sheet: "Latest Handoff Datetime" for d3763608-...md (row 6)
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H6").Value = "2016-08-24 20:43:06"
